$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 31, pushing the existing rows 31-53 down to 32-54.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Cells.Item(31, 1).Value = 2
$ws.Cells.Item(31, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44741
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = 100112026
$ws.Cells.Item(31, 7).Value = "Haba"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 700
$ws.Cells.Item(31, 11).Value = 14000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 14500
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 580
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
